$d = $word.ActiveDocument

# --- Paragraph 2: "随记: 主要原因是..." -> drop pPr rFonts hint, append "。" run ---
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5552D1FE" w14:textId="67499645" w:rsidR="00C4356A" w:rsidRDefault="00C4356A"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>随记</w:t></w:r><w:r w:rsidR="00894998"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00894998"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>主要原因是找不到工作，打算写个游戏，丰富简历，最后决定写个自传，顺便自嘲</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(2).Range.InsertXML($xml2)

# --- Paragraph 5: "建立git仓库..." -> split run, wrap "sourcetree" with proofErr ---
$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="13A3B651" w14:textId="51B7AF8F" w:rsidR="00C4356A" w:rsidRDefault="00C4356A" w:rsidP="00C4356A"><w:pPr><w:pStyle w:val="a5"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>建立git仓库，使用工具</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>sourcetree</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>管理，初步划分目录</w:t></w:r><w:r w:rsidR="00366030"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，建立对开发有利的文档：开发日记、问题记录、UI规范等。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(5).Range.InsertXML($xml5)

# --- Paragraph 9 (last paragraph in body): append the new day's note runs ---
# The last paragraph's Range includes the final body paragraph mark; replacing
# that full range with InsertXML leaves a stray empty paragraph behind because
# there's no following paragraph to "host" the mark. Exclude the trailing mark
# from the target range (End - 1) while still supplying a full <w:p> (with its
# pPr) in the injected XML so the paragraph's identity/properties are kept.
$p9 = $d.Paragraphs(9)
$r9 = $d.Range($p9.Range.Start, $p9.Range.End - 1)
$xml9 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3458AE4C" w14:textId="12ACEBF9" w:rsidR="00FB7403" w:rsidRDefault="007A6BFB" w:rsidP="0098137A"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>随记:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今天怎么也得把登录界面整好，然后淡入</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>淡出黑场</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>转场，然后设计下一场场景</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r9.InsertXML($xml9)
